$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the first visitor's e-mail address (was test6767@gmail.com)
$ws.Range("A2").Value = "SeleniumTest+v20190913120541@gmail.com"

# 2. The old value used to be a mailto: hyperlink - drop that styling, it no
#    longer points anywhere meaningful for the new generated address.
$ws.Hyperlinks.Delete()

# 3. Clear the leftover blue/underlined hyperlink font from A2 so it matches
#    the plain formatting used by the rest of the column (copy A1's format).
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. The second visitor's e-mail used to live in column B - move it back
#    under column A (was ajay.rsin@gmail.com) with the freshly generated
#    address.
$ws.Range("A3").Value = "SeleniumTest+v20190913120708@gmail.com"
$ws.Range("B3").ClearContents()
